$d = $word.ActiveDocument

# Row 6  - 1.Academics / University Result : Count {{u1}} -> 8
$d.Tables.Item(2).Rows.Item(6).Cells.Item(2).Range.Text = "8"

# Row 8  - 2.1 Journal Publications: Marks Awarded by HoD 0 -> None ; Final mark 16 -> 8
$d.Tables.Item(2).Rows.Item(8).Cells.Item(4).Range.Text = "None"
$d.Tables.Item(2).Rows.Item(8).Cells.Item(6).Range.Text = "8"

# Row 9  - Highest Impact factor: >3 : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(9).Cells.Item(4).Range.Text = "None"

# Row 10 - Highest Impact factor: 1.5 to 3 : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(10).Cells.Item(4).Range.Text = "None"

# Row 11 - Highest Impact factor: 1 to 1.4 : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(11).Cells.Item(4).Range.Text = "None"

# Row 12 - 2.2 Books authored : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(12).Cells.Item(4).Range.Text = "None"

# Row 13 - Book Editor : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(13).Cells.Item(4).Range.Text = "None"

# Row 14 - 2.3 Invited lectures presented/Chairing : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(14).Cells.Item(4).Range.Text = "None"

# Row 15 - 2.4 Conference Paper Publication (Only abroad) : Marks Awarded by HoD 6 -> None
$d.Tables.Item(2).Rows.Item(15).Cells.Item(4).Range.Text = "None"

# Row 16 - 2.5 Conference Paper Publication (India) : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(16).Cells.Item(4).Range.Text = "None"

# Row 17 - 2.6 Total Research Grant Sanctioned by Faculty : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(17).Cells.Item(4).Range.Text = "None"

# Row 18 - 2.7 Seminar Grants Received : Marks Awarded by HoD 10 -> None
$d.Tables.Item(2).Rows.Item(18).Cells.Item(4).Range.Text = "None"

# Row 19 - 2.8 Patent / Copyright Applied / Published by Faculty : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(19).Cells.Item(4).Range.Text = "None"

# Row 20 - 2.9 Consultancy /Industry project carried out by Faculty : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(20).Cells.Item(4).Range.Text = "None"

# Row 22 - 3.1 No.of Programs Attended(Seminar/workshop) by Faculty:
#          Marks Awarded by HoD 0 -> None ; Final mark 9 -> 6
$d.Tables.Item(2).Rows.Item(22).Cells.Item(4).Range.Text = "None"
$d.Tables.Item(2).Rows.Item(22).Cells.Item(6).Range.Text = "6"

# Row 23 - 3.2 No.of Skill development Programme Attended... : Marks Awarded by HoD 3 -> None
$d.Tables.Item(2).Rows.Item(23).Cells.Item(4).Range.Text = "None"

# Row 24 - 3.3 No.of NPTEL & MOOC Courses Attended by Faculty : Marks Awarded by HoD 4 -> None
$d.Tables.Item(2).Rows.Item(24).Cells.Item(4).Range.Text = "None"

# Row 25 - 3.4 Initiatives on MoU / Industrial Tie-up by Faculty : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(25).Cells.Item(4).Range.Text = "None"

# Row 26 - 3.5 Achievements, Awards, Recognition, Special Contribution : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(26).Cells.Item(4).Range.Text = "None"

# Row 27 - 3.6 No.of Conference/workshop/Hackathons organized : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(27).Cells.Item(4).Range.Text = "None"

# Row 30 - 4.1 No.of Projects Guided : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(30).Cells.Item(4).Range.Text = "None"

# Row 31 - 4.2 No.of Publications with student co-authors : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(31).Cells.Item(4).Range.Text = "None"

# Row 32 - 4.3 Mentor for reputed competitions : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(32).Cells.Item(4).Range.Text = "None"

# Row 33 - 4.4 Mentoring for student award from professional societies (or) events : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(33).Cells.Item(4).Range.Text = "None"

# Row 34 - 4.5 No.of community project completed /contribution to community based activities : Marks Awarded by HoD 0 -> None
$d.Tables.Item(2).Rows.Item(34).Cells.Item(4).Range.Text = "None"

# --- Table 3: "Faculty Cumulative Metric Score" ---

# Row 3 - Score(S): 2.Professional Activities 16 -> 8 ; 3.Research 9 -> 6 ;
#         Total Marks = Sum S x(Wi) 6.1499999999999995 -> 3.3
$d.Tables.Item(3).Rows.Item(3).Cells.Item(3).Range.Text = "8"
$d.Tables.Item(3).Rows.Item(3).Cells.Item(4).Range.Text = "6"
$d.Tables.Item(3).Rows.Item(3).Cells.Item(7).Range.Text = "3.3"

# Row 5 - Weighted Score: 2.Professional Activities 4.8 -> 2.4 ;
#         3.Research 1.3499999999999999 -> 0.8999999999999999
$d.Tables.Item(3).Rows.Item(5).Cells.Item(3).Range.Text = "2.4"
$d.Tables.Item(3).Rows.Item(5).Cells.Item(4).Range.Text = "0.8999999999999999"

Write-Output "done"
